$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Range("F3").Value = 0
$ws1.Range("F4").Value = 4738
$ws1.Range("F5").Value = 0
$ws1.Range("F6").Value = 139
$ws1.Range("F7").Value = 108
$ws1.Range("F9").Value = 0
$ws1.Range("F11").Value = 212
$ws1.Range("F12").Value = 1112
$ws1.Range("C14").Value = "南昌·漫拥动漫嘉年华Pro-追光启航（取消）"
$ws1.Range("F14").Value = 0
$ws1.Range("F15").Value = 169
$ws1.Range("F16").Value = 0
$ws1.Range("F17").Value = 136
$ws1.Range("F18").Value = 0
$ws1.Range("F19").Value = 3812
$ws1.Range("F20").Value = 6141
$ws1.Range("F24").Value = 531
$ws1.Range("F26").Value = 3937
$ws1.Range("F27").Value = 390
$ws1.Range("F29").Value = 0
$ws1.Range("F31").Value = 525
$ws1.Range("F32").Value = 0
$ws1.Range("F33").Value = 0
$ws1.Range("F34").Value = 0
$ws1.Range("F35").Value = 362
$ws1.Range("F36").Value = 156
$ws1.Range("F37").Value = 1551
$ws1.Range("F38").Value = 933
$ws1.Range("F39").Value = 39
$ws1.Range("F40").Value = 48
$ws1.Range("F41").Value = 0
$ws1.Range("F42").Value = 479
$ws1.Range("F44").Value = 73
$ws1.Range("F45").Value = 567
$ws4.Range("F2").Value = 26
$ws4.Range("F3").Value = 216
$ws4.Range("F5").Value = 0
$ws4.Range("F7").Value = 108
$ws4.Range("F8").Value = 105
$ws4.Range("F9").Value = 100
$ws4.Range("F10").Value = 0
$ws4.Range("F11").Value = 728
$ws4.Range("F13").Value = 1112
$ws4.Range("F14").Value = 97
$ws4.Range("C15").Value = "南昌·漫拥动漫嘉年华Pro-追光启航（取消）"
$ws4.Range("F15").Value = 259
$ws4.Range("F16").Value = 0
$ws4.Range("F17").Value = 78
$ws4.Range("F19").Value = 105
$ws4.Range("F20").Value = 3812
$ws4.Range("F21").Value = 6141
$ws4.Range("F23").Value = 0
$ws4.Range("F24").Value = 84
$ws4.Range("F27").Value = 0
$ws4.Range("F28").Value = 0
$ws4.Range("F29").Value = 34
$ws4.Range("F30").Value = 0
$ws4.Range("F31").Value = 567
$ws4.Range("F32").Value = 525
$ws4.Range("F33").Value = 135
$ws4.Range("F34").Value = 262
$ws4.Range("F36").Value = 0
$ws4.Range("F37").Value = 156
$ws4.Range("F38").Value = 1551
$ws4.Range("F39").Value = 0
$ws4.Range("F40").Value = 39
$ws4.Range("F41").Value = 48
$ws4.Range("F43").Value = 479
$ws4.Range("F44").Value = 477
$ws4.Range("F46").Value = 567

Write-Output "Edit complete"
